# Update cryptocurrency price (D) and 1h volume change (E) figures
# for the refreshed data pull (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that Excel would otherwise auto-detect as a number
# (e.g. "516.71") while keeping it stored as plain text, matching the source data,
# and without leaving a residual "quote prefix" text style on the cell.
function Set-TextCell($address, $text) {
    $range = $ws.Range($address)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '59.397.30'
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("D3").Value = '2.640.77'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  +0.13%  '
Set-TextCell "D5" '516.71'
$ws.Range("E5").Value = '  -0.81%  '
Set-TextCell "D6" '146.50'
$ws.Range("E6").Value = '  -1.53%  '
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("D9").Value = '2.672.40'
$ws.Range("E9").Value = '  +1.63%  '
Set-TextCell "D10" '6.46'
$ws.Range("E10").Value = '  +1.77%  '
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("E13").Value = '  -1.56%  '
$ws.Range("D14").Value = '3.104.88'
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").Value = '59.362.27'
$ws.Range("E15").Value = '  -1.32%  '
Set-TextCell "D16" '21.26'
$ws.Range("E16").Value = '  +0.13%  '
Set-TextCell "D17" '0.0000138'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '2.664.16'
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("E19").Value = '  -0.19%  '
Set-TextCell "D20" '346.57'
$ws.Range("E20").Value = '  +1.38%  '
Set-TextCell "D21" '10.51'
$ws.Range("E21").Value = '  +0.82%  '
Set-TextCell "D22" '6.17'
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("E23").Value = '  +0.15%  '
Set-TextCell "D24" '61.51'
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("D26").Value = '2.759.05'
$ws.Range("E26").Value = '  +0.89%  '
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("D29").Value = '0.0₃0819'
$ws.Range("E29").Value = '  +1.39%  '
$ws.Range("E30").Value = '  +2.10%  '
$ws.Range("E31").Value = '  -0.31%  '
Set-TextCell "D32" '6.54'
$ws.Range("E32").Value = '  +9.58%  '
Set-TextCell "D33" '19.11'
$ws.Range("E33").Value = '  +0.80%  '
$ws.Range("E34").Value = '  -0.11%  '
Set-TextCell "D35" '150.03'
$ws.Range("E35").Value = '  -0.31%  '
Set-TextCell "D36" '1.04'
$ws.Range("E36").Value = '  +13.40%  '
Set-TextCell "D37" '4.07'
$ws.Range("E37").Value = '  +3.16%  '
$ws.Range("E38").Value = '  +3.31%  '
Set-TextCell "D39" '0.868'
$ws.Range("E39").Value = '  +1.09%  '
Set-TextCell "D40" '36.75'
$ws.Range("E40").Value = '  +0.59%  '
Set-TextCell "D41" '3.73'
$ws.Range("E41").Value = '  +2.71%  '
Set-TextCell "D42" '1.42'
$ws.Range("E42").Value = '  -0.15%  '
Set-TextCell "D43" '286.03'
$ws.Range("E43").Value = '  -2.14%  '
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("E45").Value = '  -1.22%  '
Set-TextCell "D46" '0.994'
$ws.Range("E46").Value = '  -0.47%  '
Set-TextCell "D47" '19.73'
$ws.Range("E47").Value = '  +1.63%  '
Set-TextCell "D48" '0.0544'
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("E49").Value = '  +0.80%  '
Set-TextCell "D50" '4.75'
$ws.Range("E51").Value = '  -1.31%  '
